# Bwaise sysC LCA results - thorough update of results
# 1) Merge the two separate "item35"/"item36" transportation rows into a single
#    "Trucking [tonne*km]" item (two SanUnit sub-rows + one combined Total row),
#    which removes one subtotal row from the sheet.
# 2) Refresh all of the re-computed LCA numbers (construction table is untouched;
#    transportation, stream and "other" tables get refreshed values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: collapse the Transportation table ------------------------------
# Before:  row30 item35/C3 (+row31 its own Total)
#          row32 item36/C4 (+row33 its own Total)
#          row34 Sum/All
# After:   row30 Trucking/C3
#          row31 (blank)/C4
#          row32 (blank)/Total (combined)
#          row33 Sum/All
# Delete row 33 (item36's own "Total" row) - this shifts every row below it
# up by one and keeps row 30's existing "A30:A31" merge intact.
$ws.Rows.Item(33).Delete()

# Re-merge the label cell across the (now) three transportation-item rows.
$ws.Range("A30:A31").UnMerge()
$ws.Range("A30:A32").Merge()
# Merging re-derives per-cell borders for the merged block (top/middle/bottom
# slices); restore the original uniform thin-box style used everywhere else.
$ws.Range("A30:A32").Borders.LineStyle = 1

# Row 30: item35 -> Trucking, values refreshed
$ws.Range("A30").Value = "Trucking [tonne*km]"
$ws.Range("B30").Value = "C3"
$ws.Range("C30").Value = 9329531.894329509
$ws.Range("D30").Value = 0.7769153898704085
$ws.Range("E30").Value = 1809929.187499925
$ws.Range("F30").Value = 0.7769153898704086

# Row 31: was item36's own "Total" label -> now blank label + C4 SanUnit
$ws.Range("A31").Value = ""
$ws.Range("B31").Value = "C4"
$ws.Range("C31").Value = 2678895.298605489
$ws.Range("D31").Value = 0.2230846101295915
$ws.Range("E31").Value = 519705.6879294649
$ws.Range("F31").Value = 0.2230846101295916

# Row 32: was item36 label -> now blank label + combined Total
$ws.Range("A32").Value = ""
$ws.Range("B32").Value = "Total"
$ws.Range("C32").Value = 12008427.192935
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 2329634.87542939
$ws.Range("F32").Value = 1

# Row 33 (Sum/All, shifted up from old row 34) keeps its existing values.

# --- Step 2: refresh the Stream (C_* gas/liquid/solid) table ----------------
# These rows shifted up by 1 (old 37-46 -> new 36-45) via the delete above;
# only the computed C (GlobalWarming) / D (Category Ratio) numbers change.
$ws.Range("C37").Value = 49239684.13599751
$ws.Range("D37").Value = 6.683454743366449

$ws.Range("C38").Value = 4993212.770603454
$ws.Range("D38").Value = 0.6777442252504312

$ws.Range("D39").Value = -0.4339971583741242

$ws.Range("D40").Value = -4.812574815883912

$ws.Range("D41").Value = -0.1683047788592867

$ws.Range("D42").Value = -0.1555974368136683

$ws.Range("D43").Value = -0.4440753306578322

$ws.Range("D44").Value = -0.3466494480280568

$ws.Range("C45").Value = 7367399.949086143

Write-Output "done"
